$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 296.125
$ws.Range("I9").Value = 228.16667
$ws.Range("K9").Value = 228.16667
$ws.Range("M9").Value = -59.16667000000001
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = 0
$ws.Range("H53").Value = 382.7143
$ws.Range("I53").Value = 89.28570999999999
$ws.Range("J53").Value = 676.1429000000001
$ws.Range("K53").Value = 89.28570999999999
$ws.Range("L53").Value = 676.1429000000001
$ws.Range("M53").Value = 547.71429
$ws.Range("N53").Value = -1950.1429
$ws.Range("H76").Value = 8666.666999999999
$ws.Range("J76").Value = 8666.666999999999
$ws.Range("L76").Value = 8666.666999999999
$ws.Range("N76").Value = -9296.666999999999
$ws.Range("H79").Value = 8666.666999999999
$ws.Range("J79").Value = 8666.666999999999
$ws.Range("L79").Value = 8666.666999999999
$ws.Range("N79").Value = -10850.667
$ws.Range("H80").Value = 1182.2222
$ws.Range("I80").Value = 1199.8572
$ws.Range("J80").Value = 1176.05
$ws.Range("K80").Value = 3599.5716
$ws.Range("L80").Value = 3528.15
$ws.Range("M80").Value = -2601.5716
$ws.Range("N80").Value = -5524.15
$ws.Range("H83").Value = 1182.2222
$ws.Range("I83").Value = 1199.8572
$ws.Range("J83").Value = 1176.05
$ws.Range("K83").Value = 10798.7148
$ws.Range("L83").Value = 10584.45
$ws.Range("M83").Value = -5806.7148
$ws.Range("N83").Value = -20568.45
$ws.Range("H88").Value = 2850.0908
$ws.Range("I88").Value = 1725
$ws.Range("J88").Value = 3493
$ws.Range("K88").Value = 1725
$ws.Range("L88").Value = 3493
$ws.Range("M88").Value = -1319
$ws.Range("N88").Value = -4305
$ws.Range("H91").Value = 2850.0908
$ws.Range("I91").Value = 1725
$ws.Range("J91").Value = 3493
$ws.Range("K91").Value = 1725
$ws.Range("L91").Value = 3493
$ws.Range("M91").Value = -321
$ws.Range("N91").Value = -6301
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = 0
$ws.Range("H112").Value = 1676.25
$ws.Range("J112").Value = 1676.25
$ws.Range("L112").Value = 5028.75
$ws.Range("N112").Value = -7244.75
$ws.Range("H131").Value = 11840
$ws.Range("I131").Value = 10095
$ws.Range("K131").Value = 30285
$ws.Range("M131").Value = -25245
$ws.Range("H138").Value = 2700.7285
$ws.Range("I138").Value = 1342.6765
$ws.Range("J138").Value = 3983.3333
$ws.Range("K138").Value = 4028.0295
$ws.Range("L138").Value = 11949.9999
$ws.Range("M138").Value = 1111.9705
$ws.Range("N138").Value = -22229.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 2500
$ws.Range("J27").Value = 2500
$ws.Range("L27").Value = 2500
$ws.Range("N27").Value = -2868
$ws.Range("H32").Value = 19232514
$ws.Range("I32").Value = 20835098
$ws.Range("J32").Value = 1521.75
$ws.Range("K32").Value = 20835098
$ws.Range("L32").Value = 1521.75
$ws.Range("M32").Value = -20834811
$ws.Range("N32").Value = -2095.75
$ws.Range("H33").Value = 2513
$ws.Range("I33").Value = 2513
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 2513
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -2184
$ws.Range("H36").Value = 168947
$ws.Range("I36").Value = 168947
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 168947
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -168601
$ws.Range("H61").Value = 776.9722
$ws.Range("I61").Value = 763.8823
$ws.Range("K61").Value = 763.8823
$ws.Range("M61").Value = -551.8823
$ws.Range("H132").Value = 2405.238
$ws.Range("I132").Value = 2453.6
$ws.Range("J132").Value = 1438
$ws.Range("K132").Value = 7360.799999999999
$ws.Range("L132").Value = 4314
$ws.Range("M132").Value = -4830.799999999999
$ws.Range("N132").Value = -9374
$ws.Range("H136").Value = 776.9722
$ws.Range("I136").Value = 763.8823
$ws.Range("K136").Value = 2291.6469
$ws.Range("M136").Value = 258.3531000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1269.1177
$ws.Range("I20").Value = 1232.75
$ws.Range("J20").Value = 1356.4
$ws.Range("K20").Value = 1232.75
$ws.Range("L20").Value = 1356.4
$ws.Range("M20").Value = -985.75
$ws.Range("N20").Value = -1850.4
$ws.Range("H86").Value = 3191.6956
$ws.Range("I86").Value = 2779.4285
$ws.Range("J86").Value = 3833
$ws.Range("K86").Value = 2779.4285
$ws.Range("L86").Value = 3833
$ws.Range("M86").Value = -1656.4285
$ws.Range("N86").Value = -6079
$ws.Range("H89").Value = 3191.6956
$ws.Range("I89").Value = 2779.4285
$ws.Range("J89").Value = 3833
$ws.Range("K89").Value = 13897.1425
$ws.Range("L89").Value = 19165
$ws.Range("M89").Value = -8281.1425
$ws.Range("N89").Value = -30397
$ws.Range("H99").Value = 26791.555
$ws.Range("I99").Value = 29300.25
$ws.Range("K99").Value = 29300.25
$ws.Range("M99").Value = -27802.25
$ws.Range("H134").Value = 1436.909
$ws.Range("I134").Value = 1263.0625
$ws.Range("K134").Value = 3789.1875
$ws.Range("M134").Value = -1254.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3637.6
$ws.Range("I105").Value = 6445
$ws.Range("J105").Value = 1766
$ws.Range("K105").Value = 6445
$ws.Range("L105").Value = 1766
$ws.Range("M105").Value = -4698
$ws.Range("N105").Value = -5260
$ws.Range("H132").Value = 2390.484
$ws.Range("I132").Value = 2203.2
$ws.Range("K132").Value = 6609.599999999999
$ws.Range("M132").Value = -4079.599999999999
$ws.Range("H134").Value = 4560.6665
$ws.Range("I134").Value = 3355.1052
$ws.Range("J134").Value = 16013.5
$ws.Range("K134").Value = 10065.3156
$ws.Range("L134").Value = 48040.5
$ws.Range("M134").Value = -7530.3156
$ws.Range("N134").Value = -53110.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1071.591
$ws.Range("I5").Value = 731.8
$ws.Range("J5").Value = 1799.7142
$ws.Range("K5").Value = 2195.4
$ws.Range("L5").Value = 5399.142599999999
$ws.Range("M5").Value = -2083.4
$ws.Range("N5").Value = -5623.142599999999
$ws.Range("H11").Value = 3107.3157
$ws.Range("I11").Value = 4130
$ws.Range("K11").Value = 12390
$ws.Range("M11").Value = -12250
$ws.Range("H26").Value = 88.71429000000001
$ws.Range("I26").Value = 90.40000000000001
$ws.Range("J26").Value = 84.5
$ws.Range("K26").Value = 271.2
$ws.Range("L26").Value = 253.5
$ws.Range("M26").Value = 16.79999999999995
$ws.Range("N26").Value = -829.5
$ws.Range("H68").Value = 1343.1818
$ws.Range("I68").Value = 982.1429000000001
$ws.Range("K68").Value = 2946.4287
$ws.Range("M68").Value = -2135.4287
$ws.Range("H69").Value = 2800
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 4000
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = -5189
$ws.Range("N69").Value = -13622
$ws.Range("H71").Value = 1343.1818
$ws.Range("I71").Value = 982.1429000000001
$ws.Range("K71").Value = 8839.286100000001
$ws.Range("M71").Value = -4783.286100000001
$ws.Range("H72").Value = 2800
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 4000
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 36000
$ws.Range("M72").Value = -13944
$ws.Range("N72").Value = -44112
$ws.Range("H135").Value = 1071.591
$ws.Range("I135").Value = 731.8
$ws.Range("J135").Value = 1799.7142
$ws.Range("K135").Value = 6586.2
$ws.Range("L135").Value = 16197.4278
$ws.Range("M135").Value = -4051.2
$ws.Range("N135").Value = -21267.4278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 50000
$ws.Range("L42").Value = 50000
$ws.Range("N42").Value = -50970
$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -52350
$ws.Range("H122").Value = 5030.3887
$ws.Range("I122").Value = 3769.8
$ws.Range("J122").Value = 11333.333
$ws.Range("K122").Value = 11309.4
$ws.Range("L122").Value = 33999.999
$ws.Range("M122").Value = -8859.400000000001
$ws.Range("N122").Value = -38899.999
$ws.Range("H126").Value = 6640
$ws.Range("I126").Value = 7700
$ws.Range("K126").Value = 23100
$ws.Range("M126").Value = -20630
$ws.Range("H132").Value = 1631.4783
$ws.Range("I132").Value = 1523.8182
$ws.Range("K132").Value = 4571.4546
$ws.Range("M132").Value = -2041.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 872.38464
$ws.Range("I22").Value = 661.3333
$ws.Range("J22").Value = 899.913
$ws.Range("K22").Value = 661.3333
$ws.Range("L22").Value = 899.913
$ws.Range("M22").Value = -366.3333
$ws.Range("N22").Value = -1489.913
$ws.Range("H27").Value = 872.38464
$ws.Range("I27").Value = 661.3333
$ws.Range("J27").Value = 899.913
$ws.Range("K27").Value = 661.3333
$ws.Range("L27").Value = 899.913
$ws.Range("M27").Value = -554.3333
$ws.Range("N27").Value = -1113.913
$ws.Range("H93").Value = 4197.885
$ws.Range("I93").Value = 3201.7222
$ws.Range("J93").Value = 6439.25
$ws.Range("K93").Value = 3201.7222
$ws.Range("L93").Value = 6439.25
$ws.Range("M93").Value = -1953.7222
$ws.Range("N93").Value = -8935.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1065.091
$ws.Range("I107").Value = 1018.8
$ws.Range("J107").Value = 1103.6666
$ws.Range("K107").Value = 3056.4
$ws.Range("L107").Value = 3310.9998
$ws.Range("M107").Value = -1136.4
$ws.Range("N107").Value = -7150.9998
$ws.Range("H109").Value = 41095.24
$ws.Range("J109").Value = 41095.24
$ws.Range("L109").Value = 41095.24
$ws.Range("N109").Value = -43869.24
$ws.Range("H124").Value = 211499.8
$ws.Range("J124").Value = 211499.8
$ws.Range("L124").Value = 211499.8
$ws.Range("N124").Value = -221319.8
$ws.Range("H132").Value = 1429.1578
